$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of the language/value rows (rows 2-21), sorted by value
# descending, with the "Swedish" and "Uzbek" rows removed entirely.
$data = @(
    @("English", 21.33441241043863),
    @("Chinese", 19.80068753553106),
    @("Spanish", 6.302352769105785),
    @("Arabic", 4.126413723874358),
    @("German", 4.092552817610436),
    @("Japanese", 3.754450145723418),
    @("Malay-Indonesian", 3.162831628224545),
    @("Russian", 3.044977131228503),
    @("Portuguese", 2.799084056553205),
    @("French", 2.504884484349422),
    @("Turkish", 2.04975478731145),
    @("Italian", 1.870202570673083),
    @("Korean", 1.6969251959642),
    @("Dutch", 1.207422086079762),
    @("Polish", 1.001297198194399),
    @("Persian", 0.976617988020862),
    @("Urdu", 0.9281897854685975),
    @("Thai", 0.9065782757939196),
    @("Bengali", 0.819420966026593),
    @("Vietnamese", 0.8080534721724674)
)

# Clear out the old data rows (2-23) first, then remove the rows that are
# no longer needed so the sheet ends up with exactly rows 1-21.
$ws.Range("A2:B23").ClearContents()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Delete the now-unused trailing rows (22 and 23) so the sheet dimension
# shrinks to A1:B21.
$ws.Range("A22:B23").Delete()
